# Auto-considered PowerShell COM-interop script
# Actualizado 6 de septiembre 2020 - adds rows for 2020-09-03..09-06 to each dept sheet
$wb = $excel.ActiveWorkbook

# ---- Sheet "bn" ----
$ws = $wb.Worksheets.Item("bn")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Beni"
$ws.Cells.Item(179, 3).Value = 82
$ws.Cells.Item(179, 4).Value = 1
$ws.Cells.Item(179, 5).Value = 10

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Beni"
$ws.Cells.Item(180, 3).Value = 18
$ws.Cells.Item(180, 4).Value = 1
$ws.Cells.Item(180, 5).Value = 9

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Beni"
$ws.Cells.Item(181, 3).Value = 14
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 1

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Beni"
$ws.Cells.Item(182, 3).Value = 7
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 0

$ws.Range("C183:E185").Select()

# ---- Sheet "cb" ----
$ws = $wb.Worksheets.Item("cb")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B183").PasteSpecial(-4122)
$ws.Range("E178").Copy()
$ws.Range("E179:E183").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Cochabamba"
$ws.Cells.Item(179, 3).Value = 43
$ws.Cells.Item(179, 4).Value = 8
$ws.Cells.Item(179, 5).Value = 130

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Cochabamba"
$ws.Cells.Item(180, 3).Value = 30
$ws.Cells.Item(180, 4).Value = 5
$ws.Cells.Item(180, 5).Value = 90

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Cochabamba"
$ws.Cells.Item(181, 3).Value = 26
$ws.Cells.Item(181, 4).Value = 7
$ws.Cells.Item(181, 5).Value = 80

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Cochabamba"
$ws.Cells.Item(182, 3).Value = 9
$ws.Cells.Item(182, 4).Value = 5
$ws.Cells.Item(182, 5).Value = 100

$ws.Range("C184").Select()

# ---- Sheet "ch" ----
$ws = $wb.Worksheets.Item("ch")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Chuquisaca"
$ws.Cells.Item(179, 3).Value = 73
$ws.Cells.Item(179, 4).Value = 4
$ws.Cells.Item(179, 5).Value = 51

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Chuquisaca"
$ws.Cells.Item(180, 3).Value = 63
$ws.Cells.Item(180, 4).Value = 7
$ws.Cells.Item(180, 5).Value = 64

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Chuquisaca"
$ws.Cells.Item(181, 3).Value = 50
$ws.Cells.Item(181, 4).Value = 6
$ws.Cells.Item(181, 5).Value = 55

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Chuquisaca"
$ws.Cells.Item(182, 3).Value = 60
$ws.Cells.Item(182, 4).Value = 2
$ws.Cells.Item(182, 5).Value = 46

$ws.Range("C183:E185").Select()

# ---- Sheet "lp" ----
$ws = $wb.Worksheets.Item("lp")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "La Paz"
$ws.Cells.Item(179, 3).Value = 108
$ws.Cells.Item(179, 4).Value = 55
$ws.Cells.Item(179, 5).Value = 1016

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "La Paz"
$ws.Cells.Item(180, 3).Value = 214
$ws.Cells.Item(180, 4).Value = 11
$ws.Cells.Item(180, 5).Value = 1123

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "La Paz"
$ws.Cells.Item(181, 3).Value = 279
$ws.Cells.Item(181, 4).Value = 18
$ws.Cells.Item(181, 5).Value = 1056

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "La Paz"
$ws.Cells.Item(182, 3).Value = 164
$ws.Cells.Item(182, 4).Value = 6
$ws.Cells.Item(182, 5).Value = 1219

$ws.Range("C183:E185").Select()

# ---- Sheet "or" ----
$ws = $wb.Worksheets.Item("or")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Oruro"
$ws.Cells.Item(179, 3).Value = 102
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 40

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Oruro"
$ws.Cells.Item(180, 3).Value = 43
$ws.Cells.Item(180, 4).Value = 2
$ws.Cells.Item(180, 5).Value = 20

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Oruro"
$ws.Cells.Item(181, 3).Value = 19
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 9

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Oruro"
$ws.Cells.Item(182, 3).Value = 28
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 9

$ws.Range("C183:E185").Select()

# ---- Sheet "pn" ----
$ws = $wb.Worksheets.Item("pn")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Pando"
$ws.Cells.Item(179, 3).Value = 5
$ws.Cells.Item(179, 4).Value = 1
$ws.Cells.Item(179, 5).Value = 0

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Pando"
$ws.Cells.Item(180, 3).Value = 7
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 0

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Pando"
$ws.Cells.Item(181, 3).Value = 3
$ws.Cells.Item(181, 4).Value = 3
$ws.Cells.Item(181, 5).Value = 2

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Pando"
$ws.Cells.Item(182, 3).Value = 3
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 0

$ws.Range("C183:E185").Select()

# ---- Sheet "pt" ----
$ws = $wb.Worksheets.Item("pt")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Potosí"
$ws.Cells.Item(179, 3).Value = 101
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 95

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Potosí"
$ws.Cells.Item(180, 3).Value = 178
$ws.Cells.Item(180, 4).Value = 1
$ws.Cells.Item(180, 5).Value = 118

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Potosí"
$ws.Cells.Item(181, 3).Value = 68
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 63

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Potosí"
$ws.Cells.Item(182, 3).Value = 108
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 53

$ws.Range("C183:E185").Select()

# ---- Sheet "sc" ----
$ws = $wb.Worksheets.Item("sc")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Santa Cruz"
$ws.Cells.Item(179, 3).Value = 93
$ws.Cells.Item(179, 4).Value = 14
$ws.Cells.Item(179, 5).Value = 359

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Santa Cruz"
$ws.Cells.Item(180, 3).Value = 122
$ws.Cells.Item(180, 4).Value = 25
$ws.Cells.Item(180, 5).Value = 402

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Santa Cruz"
$ws.Cells.Item(181, 3).Value = 78
$ws.Cells.Item(181, 4).Value = 19
$ws.Cells.Item(181, 5).Value = 325

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Santa Cruz"
$ws.Cells.Item(182, 3).Value = 27
$ws.Cells.Item(182, 4).Value = 1593
$ws.Cells.Item(182, 5).Value = 347

$ws.Range("C183:E185").Select()

# ---- Sheet "tj" ----
$ws = $wb.Worksheets.Item("tj")
$ws.Range("A178:B178").Copy()
$ws.Range("A179:B182").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(179, 1).Value = 44077
$ws.Cells.Item(179, 2).Value = "Tarija"
$ws.Cells.Item(179, 3).Value = 246
$ws.Cells.Item(179, 4).Value = 2
$ws.Cells.Item(179, 5).Value = 109

$ws.Cells.Item(180, 1).Value = 44078
$ws.Cells.Item(180, 2).Value = "Tarija"
$ws.Cells.Item(180, 3).Value = 124
$ws.Cells.Item(180, 4).Value = 3
$ws.Cells.Item(180, 5).Value = 144

$ws.Cells.Item(181, 1).Value = 44079
$ws.Cells.Item(181, 2).Value = "Tarija"
$ws.Cells.Item(181, 3).Value = 124
$ws.Cells.Item(181, 4).Value = 2
$ws.Cells.Item(181, 5).Value = 121

$ws.Cells.Item(182, 1).Value = 44080
$ws.Cells.Item(182, 2).Value = "Tarija"
$ws.Cells.Item(182, 3).Value = 122
$ws.Cells.Item(182, 4).Value = 4
$ws.Cells.Item(182, 5).Value = 149

$ws.Range("F169").Select()

# ---- Activate the sheet that was active when the workbook was last saved ----
$wb.Worksheets.Item("tj").Activate()

$wb.Worksheets.Item("tj").Range("F169").Select()
